# Fixed query issue for C3DC phs002599
#
# The TreatmentTab query (cell B5 on Sheet1) wrapped the REPLACE(...) call
# in a redundant CONCAT(...), e.g.:
#   CONCAT(REPLACE(trt.treatment_agent, ';', ', ')) AS "Treatment Agent"
# Remove the unnecessary CONCAT() wrapper so the column reads:
#   REPLACE(trt.treatment_agent, ';', ', ') AS "Treatment Agent"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cell = $ws.Range("B5")
$query = $cell.Value2
$fixedQuery = $query.Replace( `
    "CONCAT(REPLACE(trt.treatment_agent, ';', ', ')) AS ""Treatment Agent"",", `
    "REPLACE(trt.treatment_agent, ';', ', ') AS ""Treatment Agent"",")

if ($fixedQuery -ne $query) {
    $cell.Value = $fixedQuery
}

# Reflect the reviewer leaving the cursor on the fixed cell's row (column C)
# when the workbook was last saved.
$ws.Range("C5").Select() | Out-Null
